$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create a new style combo (date format + yellow fill) for row 201, matching the new cellXfs entry ---
$ws.Range("D32").Copy($ws.Range("A201"))
$ws.Range("A201").NumberFormat = "mm-dd-yy"

# --- Append new recruitment log rows (182-208) and control-uptake rows (209-237) ---
$ws.Range("A2").Copy($ws.Range("A182"))
$ws.Range("A182").Value = 45917
$ws.Range("B182").Value = "Lausanne"
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 0

$ws.Range("A2").Copy($ws.Range("A183"))
$ws.Range("A183").Value = 45916
$ws.Range("B183").Value = "Basel"
$ws.Range("C183").Value = 1
$ws.Range("D183").Value = 0

$ws.Range("A2").Copy($ws.Range("A184"))
$ws.Range("A184").Value = 45916
$ws.Range("B184").Value = "Zuerich "
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 0

$ws.Range("A2").Copy($ws.Range("A185"))
$ws.Range("A185").Value = 45915
$ws.Range("B185").Value = "Lausanne"
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 1

$ws.Range("A2").Copy($ws.Range("A186"))
$ws.Range("A186").Value = 45915
$ws.Range("B186").Value = "Bern"
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 0

$ws.Range("A2").Copy($ws.Range("A187"))
$ws.Range("A187").Value = 45915
$ws.Range("B187").Value = "Bern"
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 1

$ws.Range("A2").Copy($ws.Range("A188"))
$ws.Range("A188").Value = 45915
$ws.Range("B188").Value = "Zuerich "
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 0

$ws.Range("A2").Copy($ws.Range("A189"))
$ws.Range("A189").Value = 45915
$ws.Range("B189").Value = "Zuerich "
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 0

$ws.Range("A2").Copy($ws.Range("A190"))
$ws.Range("A190").Value = 45911
$ws.Range("B190").Value = "Bern"
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 0

$ws.Range("A2").Copy($ws.Range("A191"))
$ws.Range("A191").Value = 45910
$ws.Range("B191").Value = "Zuerich "
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 0

$ws.Range("A2").Copy($ws.Range("A192"))
$ws.Range("A192").Value = 45909
$ws.Range("B192").Value = "Lausanne"
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 1

$ws.Range("A2").Copy($ws.Range("A193"))
$ws.Range("A193").Value = 45908
$ws.Range("B193").Value = "Geneva"
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 0

$ws.Range("A2").Copy($ws.Range("A194"))
$ws.Range("A194").Value = 45908
$ws.Range("B194").Value = "Saint Gall"
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 0

$ws.Range("A2").Copy($ws.Range("A195"))
$ws.Range("A195").Value = 45908
$ws.Range("B195").Value = "Bern"
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 1

$ws.Range("A2").Copy($ws.Range("A196"))
$ws.Range("A196").Value = 45905
$ws.Range("B196").Value = "Lausanne"
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 1

$ws.Range("A2").Copy($ws.Range("A197"))
$ws.Range("A197").Value = 45904
$ws.Range("B197").Value = "Lausanne"
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 0

$ws.Range("A2").Copy($ws.Range("A198"))
$ws.Range("A198").Value = 45904
$ws.Range("B198").Value = "Zuerich "
$ws.Range("C198").Value = 1
$ws.Range("D198").Value = 0

$ws.Range("A2").Copy($ws.Range("A199"))
$ws.Range("A199").Value = 45903
$ws.Range("B199").Value = "Lausanne"
$ws.Range("C199").Value = 1
$ws.Range("D199").Value = 1

$ws.Range("A2").Copy($ws.Range("A200"))
$ws.Range("A200").Value = 45903
$ws.Range("B200").Value = "Zuerich "
$ws.Range("C200").Value = 1
$ws.Range("D200").Value = 0

$ws.Range("A201").Value = 45902
$ws.Range("B201").Value = "Zuerich "
$ws.Range("C201").Value = 1
$ws.Range("D201").Value = 0

$ws.Range("A2").Copy($ws.Range("A202"))
$ws.Range("A202").Value = 45902
$ws.Range("B202").Value = "Bern"
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 1

$ws.Range("A2").Copy($ws.Range("A203"))
$ws.Range("A203").Value = 45902
$ws.Range("B203").Value = "Bern"
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 1

$ws.Range("A2").Copy($ws.Range("A204"))
$ws.Range("A204").Value = 45902
$ws.Range("B204").Value = "Lausanne"
$ws.Range("C204").Value = 1
$ws.Range("D204").Value = 0

$ws.Range("A2").Copy($ws.Range("A205"))
$ws.Range("A205").Value = 45902
$ws.Range("B205").Value = "Basel"
$ws.Range("C205").Value = 1
$ws.Range("D205").Value = 0

$ws.Range("A2").Copy($ws.Range("A206"))
$ws.Range("A206").Value = 45902
$ws.Range("B206").Value = "Zuerich "
$ws.Range("C206").Value = 1
$ws.Range("D206").Value = 0

$ws.Range("A2").Copy($ws.Range("A207"))
$ws.Range("A207").Value = 45898
$ws.Range("B207").Value = "Lausanne"
$ws.Range("C207").Value = 1
$ws.Range("D207").Value = 0

$ws.Range("A2").Copy($ws.Range("A208"))
$ws.Range("A208").Value = 45897
$ws.Range("B208").Value = "Lausanne"
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 1

$ws.Range("A124").Copy($ws.Range("A209"))
$ws.Range("A209").Value = 45897
$ws.Range("E209").Value = 1

$ws.Range("A124").Copy($ws.Range("A210"))
$ws.Range("A210").Value = 45898
$ws.Range("E210").Value = 1

$ws.Range("A124").Copy($ws.Range("A211"))
$ws.Range("A211").Value = 45898
$ws.Range("E211").Value = 1

$ws.Range("A82").Copy($ws.Range("A212"))
$ws.Range("A212").Value = 45904
$ws.Range("E212").Value = 1

$ws.Range("A124").Copy($ws.Range("A213"))
$ws.Range("A213").Value = 45904
$ws.Range("E213").Value = 1

$ws.Range("A124").Copy($ws.Range("A214"))
$ws.Range("A214").Value = 45905
$ws.Range("E214").Value = 1

$ws.Range("A124").Copy($ws.Range("A215"))
$ws.Range("A215").Value = 45905
$ws.Range("E215").Value = 1

$ws.Range("A124").Copy($ws.Range("A216"))
$ws.Range("A216").Value = 45905
$ws.Range("E216").Value = 1

$ws.Range("A124").Copy($ws.Range("A217"))
$ws.Range("A217").Value = 45905
$ws.Range("E217").Value = 1

$ws.Range("A124").Copy($ws.Range("A218"))
$ws.Range("A218").Value = 45905
$ws.Range("E218").Value = 1

$ws.Range("A124").Copy($ws.Range("A219"))
$ws.Range("A219").Value = 45908
$ws.Range("E219").Value = 1

$ws.Range("A124").Copy($ws.Range("A220"))
$ws.Range("A220").Value = 45908
$ws.Range("E220").Value = 1

$ws.Range("A124").Copy($ws.Range("A221"))
$ws.Range("A221").Value = 45908
$ws.Range("E221").Value = 1

$ws.Range("A124").Copy($ws.Range("A222"))
$ws.Range("A222").Value = 45909
$ws.Range("E222").Value = 1

$ws.Range("A124").Copy($ws.Range("A223"))
$ws.Range("A223").Value = 45909
$ws.Range("E223").Value = 1

$ws.Range("A124").Copy($ws.Range("A224"))
$ws.Range("A224").Value = 45909
$ws.Range("E224").Value = 1

$ws.Range("A124").Copy($ws.Range("A225"))
$ws.Range("A225").Value = 45910
$ws.Range("E225").Value = 1

$ws.Range("A124").Copy($ws.Range("A226"))
$ws.Range("A226").Value = 45910
$ws.Range("E226").Value = 1

$ws.Range("A124").Copy($ws.Range("A227"))
$ws.Range("A227").Value = 45910
$ws.Range("E227").Value = 1

$ws.Range("A124").Copy($ws.Range("A228"))
$ws.Range("A228").Value = 45910
$ws.Range("E228").Value = 1

$ws.Range("A124").Copy($ws.Range("A229"))
$ws.Range("A229").Value = 45910
$ws.Range("E229").Value = 1

$ws.Range("A124").Copy($ws.Range("A230"))
$ws.Range("A230").Value = 45911
$ws.Range("E230").Value = 1

$ws.Range("A124").Copy($ws.Range("A231"))
$ws.Range("A231").Value = 45911
$ws.Range("E231").Value = 1

$ws.Range("A82").Copy($ws.Range("A232"))
$ws.Range("A232").Value = 45915
$ws.Range("E232").Value = 1

$ws.Range("A124").Copy($ws.Range("A233"))
$ws.Range("A233").Value = 45915
$ws.Range("E233").Value = 1

$ws.Range("A124").Copy($ws.Range("A234"))
$ws.Range("A234").Value = 45916
$ws.Range("E234").Value = 1

$ws.Range("A124").Copy($ws.Range("A235"))
$ws.Range("A235").Value = 45918
$ws.Range("E235").Value = 1

$ws.Range("A124").Copy($ws.Range("A236"))
$ws.Range("A236").Value = 45918
$ws.Range("E236").Value = 1

$ws.Range("A124").Copy($ws.Range("A237"))
$ws.Range("A237").Value = 45918
$ws.Range("E237").Value = 1

# --- Update sortState to reflect the newly sorted Control tracking range ---
$rng = $ws.Range("A209:E237")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A209:A237"))
$ws.Sort.SetRange($rng)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# --- Update the active selection to match the final view state ---
$ws.Range("J229").Select()
